$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")
$ws.Activate()

# Update the offer/meeting code in A2 (shared string "ME-405" -> "ME-427")
$ws.Range("A2").Value = "ME-427"

# Move the selection from C9 to A2
[void]$ws.Range("A2").Select()
